$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 - add labels next to the existing A16/B16 styled cells
# (shared-string order must be: top, bottom, tutorial, world, map, origin)
$ws.Range("B16").Value = "top"
$ws.Range("D16").Value = "bottom"
$ws.Range("F16").Value = "tutorial"

# Row 18 - "world" data
$ws.Range("A18").Value = "world"
$ws.Range("B18").Value = -253.38
$ws.Range("C18").Value = -0.501
$ws.Range("D18").Value = 316.06
$ws.Range("E18").Value = 155.179
$ws.Range("F18").Value = -116.193
$ws.Range("G18").Value = -367.024

# Row 19 - "map" data
$ws.Range("A19").Value = "map"
$ws.Range("B19").Value = 456
$ws.Range("C19").Value = 94
$ws.Range("D19").Value = 623
$ws.Range("E19").Value = 705
$ws.Range("F19").Value = 63
$ws.Range("G19").Value = 241
$ws.Range("I19").Value = 456
$ws.Range("J19").Value = 365

# Back to row 16 - "origin" label (must come after world/map to match shared-string order)
$ws.Range("I16").Value = "origin"

# Row 22-25 - distance formulas
# F25 (and its underline styling) is set up first so the new style created for
# it only carries the font change, not the #,##0.000 number format used below.
$ws.Range("F25").Formula = "=F23/F22"
$ws.Range("F25").Font.Underline = $true

$ws.Range("F22").NumberFormat = "#,##0.000"
$ws.Range("F22").Formula = "=SQRT(SUMXMY2(B18:C18,D18:E18))"
$ws.Range("F23").NumberFormat = "#,##0.000"
$ws.Range("F23").Formula = "=SQRT(SUMXMY2(B19:C19,D19:E19))"
$ws.Range("F24").Formula = "=F22/F23"
$ws.Range("F24").ClearFormats()

# Update the selected cell to match the new active cell
$ws.Range("F27").Select() | Out-Null
